$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1970.6666
$ws.Range("J7").Value = 2306
$ws.Range("L7").Value = 2306
$ws.Range("N7").Value = -2530
$ws.Range("H14").Value = 1970.6666
$ws.Range("J14").Value = 2306
$ws.Range("L14").Value = 2306
$ws.Range("N14").Value = -2688
$ws.Range("H43").Value = 4071.2354
$ws.Range("J43").Value = 4081.1
$ws.Range("L43").Value = 4081.1
$ws.Range("N43").Value = -4219.1
$ws.Range("H113").Value = 3308.85
$ws.Range("I113").Value = 2189.4443
$ws.Range("K113").Value = 2189.4443
$ws.Range("M113").Value = 1064.5557
$ws.Range("H116").Value = 14191.7
$ws.Range("I116").Value = 7106.3335
$ws.Range("K116").Value = 7106.3335
$ws.Range("M116").Value = -3664.3335
$ws.Range("H120").Value = 97248.25
$ws.Range("J120").Value = 97248.25
$ws.Range("L120").Value = 97248.25
$ws.Range("N120").Value = -106924.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = 0
$ws.Range("H138").Value = 2131.8333
$ws.Range("I138").Value = 1276.96
$ws.Range("J138").Value = 3389
$ws.Range("K138").Value = 3830.88
$ws.Range("L138").Value = 10167
$ws.Range("M138").Value = 1309.12
$ws.Range("N138").Value = -20447

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22074.154
$ws.Range("I32").Value = 22535.19
$ws.Range("J32").Value = 17187.2
$ws.Range("K32").Value = 22535.19
$ws.Range("L32").Value = 17187.2
$ws.Range("M32").Value = -22248.19
$ws.Range("N32").Value = -17761.2
$ws.Range("H61").Value = 1754.1364
$ws.Range("I61").Value = 1281.1052
$ws.Range("K61").Value = 1281.1052
$ws.Range("M61").Value = -1069.1052
$ws.Range("H74").Value = 6000012
$ws.Range("I74").Value = 6000012
$ws.Range("K74").Value = 6000012
$ws.Range("M74").Value = -5999138
$ws.Range("H77").Value = 6000012
$ws.Range("I77").Value = 6000012
$ws.Range("K77").Value = 30000060
$ws.Range("M77").Value = -29995692
$ws.Range("H102").Value = 3052.465
$ws.Range("I102").Value = 2351.5881
$ws.Range("J102").Value = 5700.222
$ws.Range("K102").Value = 2351.5881
$ws.Range("L102").Value = 5700.222
$ws.Range("M102").Value = -729.5880999999999
$ws.Range("N102").Value = -8944.222
$ws.Range("H132").Value = 1293.5122
$ws.Range("I132").Value = 1078
$ws.Range("J132").Value = 1630.25
$ws.Range("K132").Value = 3234
$ws.Range("L132").Value = 4890.75
$ws.Range("M132").Value = -704
$ws.Range("N132").Value = -9950.75
$ws.Range("H136").Value = 1754.1364
$ws.Range("I136").Value = 1281.1052
$ws.Range("K136").Value = 3843.3156
$ws.Range("M136").Value = -1293.3156

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2466.818
$ws.Range("I86").Value = 1761.75
$ws.Range("J86").Value = 4347
$ws.Range("K86").Value = 1761.75
$ws.Range("L86").Value = 4347
$ws.Range("M86").Value = -638.75
$ws.Range("N86").Value = -6593
$ws.Range("H89").Value = 2466.818
$ws.Range("I89").Value = 1761.75
$ws.Range("J89").Value = 4347
$ws.Range("K89").Value = 8808.75
$ws.Range("L89").Value = 21735
$ws.Range("M89").Value = -3192.75
$ws.Range("N89").Value = -32967
$ws.Range("H105").Value = 2406.2856
$ws.Range("I105").Value = 1268.8
$ws.Range("K105").Value = 1268.8
$ws.Range("M105").Value = 478.2
$ws.Range("H107").Value = 16698.486
$ws.Range("I107").Value = 18456.42
$ws.Range("J107").Value = 3074.5
$ws.Range("K107").Value = 18456.42
$ws.Range("L107").Value = 3074.5
$ws.Range("M107").Value = -16536.42
$ws.Range("N107").Value = -6914.5
$ws.Range("H117").Value = 68742
$ws.Range("J117").Value = 68742
$ws.Range("L117").Value = 68742
$ws.Range("N117").Value = -77920
$ws.Range("H134").Value = 2392.4
$ws.Range("I134").Value = 2357
$ws.Range("K134").Value = 7071
$ws.Range("M134").Value = -4536

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3768.6667
$ws.Range("I10").Value = 653.5
$ws.Range("J10").Value = 9999
$ws.Range("K10").Value = 653.5
$ws.Range("L10").Value = 9999
$ws.Range("M10").Value = -514.5
$ws.Range("N10").Value = -10277
$ws.Range("H31").Value = 2633931
$ws.Range("I31").Value = 3705986.8
$ws.Range("J31").Value = 2521.6365
$ws.Range("K31").Value = 3705986.8
$ws.Range("L31").Value = 2521.6365
$ws.Range("M31").Value = -3705691.8
$ws.Range("N31").Value = -3111.6365
$ws.Range("H34").Value = 2633931
$ws.Range("I34").Value = 3705986.8
$ws.Range("J34").Value = 2521.6365
$ws.Range("K34").Value = 3705986.8
$ws.Range("L34").Value = 2521.6365
$ws.Range("M34").Value = -3705784.8
$ws.Range("N34").Value = -2925.6365
$ws.Range("H62").Value = 7645.3335
$ws.Range("I62").Value = 4466.6665
$ws.Range("K62").Value = 4466.6665
$ws.Range("M62").Value = -3842.6665
$ws.Range("H65").Value = 7645.3335
$ws.Range("I65").Value = 4466.6665
$ws.Range("K65").Value = 22333.3325
$ws.Range("M65").Value = -19213.3325
$ws.Range("H105").Value = 1615.5834
$ws.Range("I105").Value = 1131.8889
$ws.Range("J105").Value = 3066.6667
$ws.Range("K105").Value = 1131.8889
$ws.Range("L105").Value = 3066.6667
$ws.Range("M105").Value = 615.1111000000001
$ws.Range("N105").Value = -6560.6667
$ws.Range("H107").Value = 561.44116
$ws.Range("I107").Value = 401.20834
$ws.Range("K107").Value = 401.20834
$ws.Range("M107").Value = 1518.79166

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 306.375
$ws.Range("I14").Value = 306.375
$ws.Range("K14").Value = 919.125
$ws.Range("M14").Value = -746.125
$ws.Range("H60").Value = 4104.72
$ws.Range("I60").Value = 204.5
$ws.Range("J60").Value = 4847.619
$ws.Range("K60").Value = 613.5
$ws.Range("L60").Value = 14542.857
$ws.Range("M60").Value = -362.5
$ws.Range("N60").Value = -15044.857

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 344.83334
$ws.Range("J2").Value = 800
$ws.Range("L2").Value = 800
$ws.Range("N2").Value = -1026
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5232
$ws.Range("H102").Value = 17468.424
$ws.Range("I102").Value = 18715.434
$ws.Range("K102").Value = 18715.434
$ws.Range("M102").Value = -17093.434
$ws.Range("H107").Value = 262.42856
$ws.Range("I107").Value = 84.666664
$ws.Range("K107").Value = 84.666664
$ws.Range("M107").Value = 1835.333336
$ws.Range("H113").Value = 1377.5
$ws.Range("I113").Value = 1370.3334
$ws.Range("J113").Value = 1399
$ws.Range("K113").Value = 1370.3334
$ws.Range("L113").Value = 1399
$ws.Range("M113").Value = 799.6666
$ws.Range("N113").Value = -5739
$ws.Range("H133").Value = 93415.89
$ws.Range("J133").Value = 102890
$ws.Range("L133").Value = 102890
$ws.Range("N133").Value = -113010

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 12933.75
$ws.Range("J14").Value = 12933.75
$ws.Range("L14").Value = 12933.75
$ws.Range("N14").Value = -13277.75
$ws.Range("H46").Value = 4584.6816
$ws.Range("I46").Value = 2003.1
$ws.Range("J46").Value = 6736
$ws.Range("K46").Value = 2003.1
$ws.Range("L46").Value = 6736
$ws.Range("M46").Value = -1815.1
$ws.Range("N46").Value = -7112
$ws.Range("H136").Value = 8476.718000000001
$ws.Range("I136").Value = 8405.333000000001
$ws.Range("K136").Value = 25215.999
$ws.Range("M136").Value = -22665.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3090.6365
$ws.Range("I81").Value = 3547.5
$ws.Range("J81").Value = 1872.3334
$ws.Range("K81").Value = 7095
$ws.Range("L81").Value = 3744.6668
$ws.Range("M81").Value = -6034
$ws.Range("N81").Value = -5866.6668
$ws.Range("H84").Value = 3090.6365
$ws.Range("I84").Value = 3547.5
$ws.Range("J84").Value = 1872.3334
$ws.Range("K84").Value = 35475
$ws.Range("L84").Value = 18723.334
$ws.Range("M84").Value = -30171
$ws.Range("N84").Value = -29331.334
$ws.Range("H116").Value = 71560
$ws.Range("J116").Value = 71560
$ws.Range("L116").Value = 71560
$ws.Range("N116").Value = -80738
$ws.Range("H122").Value = 107403.59
$ws.Range("I122").Value = 111565.664
$ws.Range("K122").Value = 334696.992
$ws.Range("M122").Value = -332246.992
$ws.Range("H132").Value = 6151.1177
$ws.Range("I132").Value = 7314.75
$ws.Range("K132").Value = 21944.25
$ws.Range("M132").Value = -19414.25
$ws.Range("H136").Value = 14919.53
$ws.Range("I136").Value = 17884.875
$ws.Range("K136").Value = 53654.625
$ws.Range("M136").Value = -51104.625
